$d = $word.ActiveDocument

$wNS = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

# --- Step 1: split the "timer #(.NTICKS(10))" paragraph into three
#     paragraphs that declare TIMER_PERIOD_NS and CLOCK_PERIOD_NS as
#     module parameters. ---
$instParaXml = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
<w:body>
<w:p>
  <w:pPr>
    <w:rPr>
      <w:rFonts w:ascii="Consolas" w:hAnsi="Consolas"/>
    </w:rPr>
  </w:pPr>
  <w:proofErr w:type="gramStart"/>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="Consolas" w:hAnsi="Consolas"/>
    </w:rPr>
    <w:t>timer</w:t>
  </w:r>
  <w:proofErr w:type="gramEnd"/>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="Consolas" w:hAnsi="Consolas"/>
    </w:rPr>
    <w:t xml:space="preserve"> </w:t>
  </w:r>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="Consolas" w:hAnsi="Consolas"/>
    </w:rPr>
    <w:t>#(</w:t>
  </w:r>
</w:p>
<w:p>
  <w:pPr>
    <w:rPr>
      <w:rFonts w:ascii="Consolas" w:hAnsi="Consolas"/>
    </w:rPr>
  </w:pPr>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="Consolas" w:hAnsi="Consolas"/>
    </w:rPr>
    <w:t xml:space="preserve">    </w:t>
  </w:r>
  <w:proofErr w:type="gramStart"/>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="Consolas" w:hAnsi="Consolas"/>
    </w:rPr>
    <w:t>.</w:t>
  </w:r>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="Consolas" w:hAnsi="Consolas"/>
    </w:rPr>
    <w:t>TIMER</w:t>
  </w:r>
  <w:proofErr w:type="gramEnd"/>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="Consolas" w:hAnsi="Consolas"/>
    </w:rPr>
    <w:t>_PERIOD_NS</w:t>
  </w:r>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="Consolas" w:hAnsi="Consolas"/>
    </w:rPr>
    <w:t>(</w:t>
  </w:r>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="Consolas" w:hAnsi="Consolas"/>
    </w:rPr>
    <w:t>80</w:t>
  </w:r>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="Consolas" w:hAnsi="Consolas"/>
    </w:rPr>
    <w:t>,</w:t>
  </w:r>
</w:p>
<w:p>
  <w:pPr>
    <w:rPr>
      <w:rFonts w:ascii="Consolas" w:hAnsi="Consolas"/>
    </w:rPr>
  </w:pPr>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="Consolas" w:hAnsi="Consolas"/>
    </w:rPr>
    <w:t xml:space="preserve">    </w:t>
  </w:r>
  <w:proofErr w:type="gramStart"/>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="Consolas" w:hAnsi="Consolas"/>
    </w:rPr>
    <w:t>.CLOCK</w:t>
  </w:r>
  <w:proofErr w:type="gramEnd"/>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="Consolas" w:hAnsi="Consolas"/>
    </w:rPr>
    <w:t>_PERIOD_NS(8))</w:t>
  </w:r>
</w:p>
</w:body>
</w:document>
'@

$instPara = $d.Paragraphs(5).Range
$instPara.InsertXML($instParaXml)

# --- Step 2: rewrite the NTICKS parameter description paragraph,
#     inserting two new paragraphs (TIMER_PERIOD_NS / CLOCK_PERIOD_NS)
#     ahead of it, each followed by a blank paragraph, and updating
#     the NTICKS text itself. ---
$nticksParaXml = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
<w:body>
<w:p>
  <w:r>
    <w:t>TIMER_PERIOD_NS is the amount of time taken for the timer to expire once armed. Units are in nanoseconds. Default value is 80ns.</w:t>
  </w:r>
</w:p>
<w:p/>
<w:p>
  <w:r>
    <w:t>CLOCK_PERIOD_NS is the period of the common clock. Units are in nanoseconds. Default value is 8ns (125MHz).</w:t>
  </w:r>
</w:p>
<w:p/>
<w:p>
  <w:r>
    <w:t>NTICKS</w:t>
  </w:r>
  <w:r>
    <w:t xml:space="preserve"> is the </w:t>
  </w:r>
  <w:r>
    <w:t xml:space="preserve">number of </w:t>
  </w:r>
  <w:r>
    <w:t xml:space="preserve">clock ticks to count until the timer expires. </w:t>
  </w:r>
  <w:r>
    <w:t>Default value is TIMER_PERIOD_NS / CLOCK_PERIOD_NS. This value can be specified directly in case the clock period is not an integral number of nanoseconds.</w:t>
  </w:r>
</w:p>
</w:body>
</w:document>
'@

$nticksPara = $d.Paragraphs(16).Range
$nticksPara.InsertXML($nticksParaXml)
